# [ADD] Proceso 1 done
# The "id" column (A) used to start at 1; this process adds 100 to every
# row so the sequence now runs 101-150. Column A3:A51 already carries the
# incremental formula (=A2+1 / shared =A{n-1}+1), so updating the seed
# value in A2 ripples through the rest of the column automatically on
# recalculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 101

# Reflect the author's final cursor position (cell A3 selected) in the
# saved sheet view.
$ws.Range("A3").Select()
